# Scaling template update. NFR and GNFR entries are required.
#
# Fills in the previously-blank GNFR_code / NFR_code cells on the "Scaling"
# sheet (columns E and F of tbl_I_scaling) so the calculated helper columns
# (GNFR_label, NFR_label, the GNFR-NFR range lookup and the GNFR/NFR
# consistency check) resolve instead of erroring with #N/A.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Scaling")

# GNFR_code (column E) was blank for these rows -> fill with the matching
# GNFR code so the row's NFR_code can be validated against it.
$ws.Range("E2").Value = "B"
$ws.Range("E3").Value = "B"
$ws.Range("E4").Value = "B"
$ws.Range("E5").Value = "B"
$ws.Range("E7").Value = "K"
$ws.Range("E8").Value = "K"
$ws.Range("E9").Value = "J"
$ws.Range("E10").Value = "E"

# NFR_code (column F) was blank for these rows (only GNFR_code "*" / "B"
# was set) -> fill with the wildcard "*" value.
$ws.Range("F6").Value = "*"
$ws.Range("F11").Value = "*"

$wb.Application.Calculate()

# Restore the cursor to the cell it ends up on after the edits.
$ws.Range("F14").Select()
